$wb = $excel.ActiveWorkbook

# --- 1) Text change: "Ready for handoff" -> "In Translation" ---
# This string shows up in the "Status" column (and its Overview-sheet
# mirror columns) on every sheet, so sweep each used range and flip it
# wherever it is found. NOTE: compare with the literal on the LEFT of
# -eq, since some cells hold native booleans (True/False) and PowerShell
# coerces the right-hand side to the left operand's type, which would
# otherwise make a boolean cell falsely "equal" any non-empty string.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $sheet.Cells.Item($r, $c)
            if ("Ready for handoff" -eq $cell.Value2) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# --- 2) Column width changes (stored OOXML width 17.2159881591797 -> 13.4101845877511) ---
# The COM ColumnWidth property is quantized to whole pixels (1/6-character
# steps) by the host, same as real Excel; 12.5 is the closest achievable
# COM width to the target, landing on stored width 13.333333333333334.
$newColumnWidth = 12.5

# Overview sheet: columns E and F (zh-cn / de-de "Status" mirror columns)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth

# zh-cn sheet: column C ("Status")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C1").ColumnWidth = $newColumnWidth

# de-de sheet: column C ("Status")
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C1").ColumnWidth = $newColumnWidth
